$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Duel Decks Jace vs. Chandra Tokens (TDD2)'
$ws.Cells.Item(2, 1).Value = 'Aethersnipe'
$ws.Cells.Item(3, 1).Value = 'Air Elemental'
$ws.Cells.Item(4, 1).Value = 'Ancestral Vision'
$ws.Cells.Item(5, 1).Value = 'Bottle Gnomes'
$ws.Cells.Item(6, 1).Value = 'Brine Elemental'
$ws.Cells.Item(7, 1).Value = 'Chandra Nalaar'
$ws.Cells.Item(8, 1).Value = 'Chartooth Cougar'
$ws.Cells.Item(9, 1).Value = 'Condescend'
$ws.Cells.Item(10, 1).Value = 'Cone of Flame'
$ws.Cells.Item(11, 1).Value = 'Counterspell'
$ws.Cells.Item(12, 1).Value = 'Daze'
$ws.Cells.Item(13, 1).Value = 'Demonfire'
$ws.Cells.Item(14, 1).Value = 'Errant Ephemeron'
$ws.Cells.Item(15, 1).Value = 'Fact or Fiction'
$ws.Cells.Item(16, 1).Value = 'Fathom Seer'
$ws.Cells.Item(17, 1).Value = 'Fireball'
$ws.Cells.Item(18, 1).Value = 'Fireblast'
$ws.Cells.Item(19, 1).Value = 'Firebolt'
$ws.Cells.Item(20, 1).Value = 'Fireslinger'
$ws.Cells.Item(21, 1).Value = 'Flame Javelin'
$ws.Cells.Item(22, 1).Value = 'Flamekin Brawler'
$ws.Cells.Item(23, 1).Value = 'Flametongue Kavu'
$ws.Cells.Item(24, 1).Value = 'Flamewave Invoker'
$ws.Cells.Item(25, 1).Value = 'Fledgling Mawcor'
$ws.Cells.Item(26, 1).Value = 'Furnace Whelp'
$ws.Cells.Item(27, 1).Value = 'Guile'
$ws.Cells.Item(28, 1).Value = 'Gush'
$ws.Cells.Item(29, 1).Value = 'Hostility'
$ws.Cells.Item(30, 1).Value = 'Incinerate'
$ws.Cells.Item(31, 1).Value = 'Ingot Chewer'
$ws.Cells.Item(32, 1).Value = 'Inner-Flame Acolyte'
$ws.Cells.Item(33, 1).Value = 'Island'
$ws.Cells.Item(34, 1).Value = 'Island'
$ws.Cells.Item(35, 1).Value = 'Island'
$ws.Cells.Item(36, 1).Value = 'Island'
$ws.Cells.Item(37, 1).Value = 'Jace Beleren'
$ws.Cells.Item(38, 1).Value = 'Keldon Megaliths'
$ws.Cells.Item(39, 1).Value = 'Magma Jet'
$ws.Cells.Item(40, 1).Value = 'Man-o''-War'
$ws.Cells.Item(41, 1).Value = 'Martyr of Frost'
$ws.Cells.Item(42, 1).Value = 'Mind Stone'
$ws.Cells.Item(43, 1).Value = 'Mountain'
$ws.Cells.Item(44, 1).Value = 'Mountain'
$ws.Cells.Item(45, 1).Value = 'Mountain'
$ws.Cells.Item(46, 1).Value = 'Mountain'
$ws.Cells.Item(47, 1).Value = 'Mulldrifter'
$ws.Cells.Item(48, 1).Value = 'Ophidian'
$ws.Cells.Item(49, 1).Value = 'Oxidda Golem'
$ws.Cells.Item(50, 1).Value = 'Pyre Charger'
$ws.Cells.Item(51, 1).Value = 'Quicksilver Dragon'
$ws.Cells.Item(52, 1).Value = 'Rakdos Pit Dragon'
$ws.Cells.Item(53, 1).Value = 'Repulse'
$ws.Cells.Item(54, 1).Value = 'Riftwing Cloudskate'
$ws.Cells.Item(55, 1).Value = 'Seal of Fire'
$ws.Cells.Item(56, 1).Value = 'Slith Firewalker'
$ws.Cells.Item(57, 1).Value = 'Soulbright Flamekin'
$ws.Cells.Item(58, 1).Value = 'Spire Golem'
$ws.Cells.Item(59, 1).Value = 'Terrain Generator'
$ws.Cells.Item(60, 1).Value = 'Voidmage Apprentice'
$ws.Cells.Item(61, 1).Value = 'Wall of Deceit'
$ws.Cells.Item(62, 1).Value = 'Waterspout Djinn'
$ws.Cells.Item(63, 1).Value = 'Willbender'
